$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 24; $row -le 81; $row++) {
    $suffix = 3997 + ($row - 24)
    $ws.Cells.Item($row, 19).Value = "https://orcid.org/0000-0003-2195-$suffix"
}
